$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole table (previously in columns B:C, rows 3-8) moves one column to
# the left and up so the title sits on row 1 and the table starts on row 3.
# Deleting the empty leading column/rows shifts the existing content (and
# keeps its formatting / column widths) rather than re-typing everything.
$ws.Columns.Item(1).Delete()
$ws.Range("A1:A2").EntireRow.Delete()

# Update selection to match the saved view state
$ws.Range("A3").Select()
